# Daily update at 8 AM UTC: append the next day's row of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 89
$prevRow = 88

# New day's values
$ws.Cells.Item($newRow, 1).Value = 46038
$ws.Cells.Item($newRow, 2).Value = 204
$ws.Cells.Item($newRow, 3).Value = 213
$ws.Cells.Item($newRow, 4).Value = 200

# Match the formatting of the cell above (date style, etc.) like Excel does
# when you fill/append a new row following the existing pattern.
$ws.Range("A$prevRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0
